# Updates the cryptos list Price (column D) and Volume(1h) (column E) values
# to the latest snapshot scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = '<new price text>'; E = '<new volume text>' }
$updates = @{
    2 = @{ D='26.093.47'; E='  -1.73%  ' }
    3 = @{ D='1.666.25'; E='  -1.20%  ' }
    4 = @{ E='  -0.04%  ' }
    5 = @{ D='216.69'; E='  +0.05%  ' }
    6 = @{ D='0.5102'; E='  +2.91%  ' }
    7 = @{ E='  -0.01%  ' }
    8 = @{ D='0.2635'; E='  +1.27%  ' }
    9 = @{ E='  +4.89%  ' }
    10 = @{ D='21.51'; E='  -0.37%  ' }
    11 = @{ D='0.07394'; E='  +1.76%  ' }
    12 = @{ D='1.675.22'; E='  -0.76%  ' }
    13 = @{ D='4.505'; E='  +1.77%  ' }
    14 = @{ D='0.5790'; E='  +1.47%  ' }
    15 = @{ D='0.000008555'; E='  +3.78%  ' }
    16 = @{ D='64.17'; E='  -0.48%  ' }
    17 = @{ D='26.146.47' }
    18 = @{ D='4.918'; E='  -1.50%  ' }
    19 = @{ D='1.005'; E='  -0.04%  ' }
    20 = @{ D='10.80'; E='  +1.21%  ' }
    21 = @{ D='189.25'; E='  +4.03%  ' }
    22 = @{ D='6.203'; E='  +0.78%  ' }
    23 = @{ D='1.006'; E='  +0.02%  ' }
    24 = @{ D='144.98'; E='  +0.36%  ' }
    25 = @{ E='  +0.88%  ' }
    26 = @{ D='0.1195'; E='  +6.10%  ' }
    27 = @{ D='15.57'; E='  +2.23%  ' }
    28 = @{ D='0.06355'; E='  +14.11%  ' }
    29 = @{ E='  -1.22%  ' }
    30 = @{ D='1.315'; E='  -0.60%  ' }
    31 = @{ D='3.523'; E='  +1.62%  ' }
    32 = @{ D='3.506'; E='  +1.47%  ' }
    33 = @{ D='1.632'; E='  -0.56%  ' }
    34 = @{ D='1.013'; E='  +0.86%  ' }
    35 = @{ D='0.6082'; E='  +3.95%  ' }
    36 = @{ D='2.364'; E='  -0.68%  ' }
    37 = @{ D='2.648'; E='  +0.49%  ' }
    38 = @{ D='6.149'; E='  +4.33%  ' }
    39 = @{ D='0.01610'; E='  +1.82%  ' }
    40 = @{ D='1.074.61'; E='  +0.47%  ' }
    41 = @{ D='0.8599'; E='  +1.29%  ' }
    43 = @{ D='101.08'; E='  +3.13%  ' }
    44 = @{ D='1.814.46'; E='  -1.66%  ' }
    45 = @{ D='0.00000000114'; E='  +8.16%  ' }
    46 = @{ D='56.12'; E='  +0.08%  ' }
    47 = @{ D='1.006'; E='  +0.05%  ' }
    48 = @{ D='8.051'; E='  +0.24%  ' }
    49 = @{ E='  -0.29%  ' }
    50 = @{ D='0.4290'; E='  -0.96%  ' }
    51 = @{ D='5.922'; E='  +6.74%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('D')) {
        $cell = $ws.Cells.Item($row, 4)
        # Force text so values like trailing-zero decimals (e.g. 0.5790)
        # and thousand-dotted prices (e.g. 26.093.47) aren't coerced to numbers,
        # then drop the temporary number-format override so the cell keeps
        # its original (default) style, same as every other data cell.
        $cell.NumberFormat = '@'
        $cell.Value = $vals['D']
        $cell.ClearFormats()
    }
    if ($vals.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $vals['E']
    }
}
